$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: mark the "本地缓存localstorage" TODO as done (status + resolved date) ---
# Copy formatting from an existing "status" (C) / "resolved date" (D) cell so the
# new cells reuse the same cell styles instead of creating new ones.
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("C6").Value = "√"

$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("D6").Value = 42806                # 2017-03-12

# --- Row 21: status cell keeps its "待定" text, just restyle to match other status cells ---
$ws.Range("C2").Copy()
$ws.Range("C21").PasteSpecial(-4122)         # xlPasteFormats (value/text untouched)

# --- Row 23: mark the "没有必要性" TODO as done (status + resolved date) ---
$ws.Range("C2").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "√"

$ws.Range("D2").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 42805               # 2017-03-11

# --- Row 24: bring the new status/resolved-date columns into this row, left blank ---
$ws.Range("C2").Copy()
$ws.Range("C24").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D24").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Update the view: scroll back up and move the active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E8").Select()
